# Savesheet para diferentes bimestres
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("B4").Value = "MEC-3B-M. A. Comp; Cad / CAM"
$ws.Range("E4").Value = "-"

# Row 6
$ws.Range("B6").Value = "MEC-3B-M. A. Comp; Cad / CAM"
$ws.Range("E6").Value = "-"

# Row 7
$ws.Range("B7").Value = "MEC-3B-M. A. Comp; Cad / CAM"
$ws.Range("E7").Value = "-"

# Row 8
$ws.Range("B8").Value = "MEC-3B-M. A. Comp; Cad / CAM"
$ws.Range("E8").Value = "-"

# Row 10
$ws.Range("D10").Value = "-"
$ws.Range("E10").Value = "MEC-1A-Des. Tec. Mec."
$ws.Range("F10").Value = "MEC-1A-Des. Tec. Mec."

# Row 11
$ws.Range("D11").Value = "-"
$ws.Range("E11").Value = "MEC-1A-Des. Tec. Mec."
$ws.Range("F11").Value = "MEC-1A-Des. Tec. Mec."

# Row 12
$ws.Range("D12").Value = "-"
$ws.Range("E12").Value = "MEC-1A-Des. Tec. Mec."

# Row 14
$ws.Range("D14").Value = "-"

# Row 15
$ws.Range("D15").Value = "-"

# Row 16
$ws.Range("D16").Value = "-"
$ws.Range("E16").Value = "MEC-1A-Des. Tec. Mec."
